$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("M2").Value = 428
$ws.Range("M3").Value = 471
$ws.Range("M4").Value = 134
$ws.Range("M5").Value = 29
$ws.Range("M6").Value = 363
$ws.Range("M7").Value = 1425

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("M2").Value = 27
$ws.Range("M3").Value = 33
$ws.Range("M5").Value = 3
$ws.Range("M7").Value = 95

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("M3").Value = 14
$ws.Range("M7").Value = 31

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("M2").Value = 12
$ws.Range("M5").Value = 3
$ws.Range("M7").Value = 52

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("M3").Value = 7
$ws.Range("M6").Value = 6
$ws.Range("M7").Value = 19

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("M6").Value = 17
$ws.Range("M7").Value = 66

$ws = $wb.Worksheets.Item("New City")
$ws.Range("M3").Value = 12
$ws.Range("M6").Value = 7
$ws.Range("M7").Value = 30

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("M2").Value = 12
$ws.Range("M7").Value = 42
$ws.Range("M8").Value = 95
$ws.Range("M9").Value = 12
$ws.Range("M19").Value = 50
$ws.Range("M20").Value = 48
$ws.Range("M27").Value = 20
$ws.Range("M29").Value = 71
$ws.Range("M33").Value = 52
$ws.Range("M34").Value = 9
$ws.Range("M37").Value = 66
$ws.Range("M42").Value = 48
$ws.Range("M48").Value = 17
$ws.Range("M60").Value = 12
$ws.Range("M63").Value = 5
$ws.Range("M65").Value = 30
$ws.Range("M68").Value = 4
$ws.Range("M71").Value = 9
$ws.Range("M72").Value = 10
$ws.Range("M76").Value = 17
$ws.Range("M77").Value = 15
$ws.Range("M78").Value = 22
$ws.Range("M79").Value = 34
$ws.Range("M83").Value = 31
$ws.Range("M85").Value = 69
$ws.Range("M87").Value = 4
$ws.Range("M88").Value = 17
$ws.Range("M89").Value = 21
$ws.Range("M90").Value = 14
$ws.Range("M91").Value = 20
$ws.Range("M94").Value = 19
$ws.Range("M95").Value = 19
$ws.Range("M96").Value = 13
$ws.Range("M101").Value = 1425

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("M2").Value = 23
$ws.Range("M3").Value = 23
$ws.Range("M4").Value = 6
$ws.Range("M6").Value = 18
$ws.Range("M7").Value = 71

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("M3").Value = 2
$ws.Range("M7").Value = 17

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("M2").Value = 15
$ws.Range("M7").Value = 50

$ws = $wb.Worksheets.Item("River North")
$ws.Range("M3").Value = 4
$ws.Range("M7").Value = 17

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("M6").Value = 15
$ws.Range("M7").Value = 48

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("M3").Value = 11
$ws.Range("M7").Value = 22

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("M6").Value = 5
$ws.Range("M7").Value = 13

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("M5").Value = 1
$ws.Range("M7").Value = 20

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("M2").Value = 12
$ws.Range("M6").Value = 8
$ws.Range("M7").Value = 34

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("M2").Value = 20
$ws.Range("M3").Value = 11
$ws.Range("M4").Value = 6
$ws.Range("M7").Value = 48

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("M3").Value = 15
$ws.Range("M7").Value = 42

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("M2").Value = 1
$ws.Range("M7").Value = 9

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("M3").Value = 5
$ws.Range("M7").Value = 19

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("M6").Value = 6
$ws.Range("M7").Value = 12

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("M4").Value = 5
$ws.Range("M7").Value = 12

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("M2").Value = 4
$ws.Range("M6").Value = 9
$ws.Range("M7").Value = 17

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("M4").Value = 5
$ws.Range("M7").Value = 21

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("M4").Value = 6
$ws.Range("M5").Value = 3
$ws.Range("M6").Value = 20

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("M6").Value = 3
$ws.Range("M7").Value = 14

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("M5").Value = 1
$ws.Range("M6").Value = 4

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("M2").Value = 5
$ws.Range("M7").Value = 12

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("M4").Value = 3
$ws.Range("M7").Value = 69

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("L4").Value = 1
$ws.Range("L7").Value = 9

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("M2").Value = 3
$ws.Range("M4").Value = 4
$ws.Range("M6").Value = 10

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("M3").Value = 5
$ws.Range("M7").Value = 15

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("K6").Value = 1
$ws.Range("K7").Value = 4
